$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 465.66666
$ws.Range("J12").Value = 449
$ws.Range("L12").Value = 449
$ws.Range("N12").Value = -789

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20924

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 19050
$ws.Range("I54").Value = 19050
$ws.Range("K54").Value = 19050
$ws.Range("M54").Value = -18564

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1620
$ws.Range("I62").Value = 1250
$ws.Range("J62").Value = 1866.6666
$ws.Range("K62").Value = 1250
$ws.Range("L62").Value = 1866.6666
$ws.Range("M62").Value = -626
$ws.Range("N62").Value = -3114.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 1620
$ws.Range("I65").Value = 1250
$ws.Range("J65").Value = 1866.6666
$ws.Range("K65").Value = 6250
$ws.Range("L65").Value = 9333.333000000001
$ws.Range("M65").Value = -3130
$ws.Range("N65").Value = -15573.333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 642
$ws.Range("I96").Value = 670.4
$ws.Range("J96").Value = 500
$ws.Range("K96").Value = 2011.2
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = -638.1999999999998
$ws.Range("N96").Value = -4246

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6500.5
$ws.Range("I116").Value = 6000
$ws.Range("K116").Value = 6000
$ws.Range("M116").Value = -2558

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2980.3333
$ws.Range("I137").Value = 2979.6667
$ws.Range("K137").Value = 8939.000100000001
$ws.Range("M137").Value = -6389.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4149.75
$ws.Range("I32").Value = 2949.75
$ws.Range("K32").Value = 2949.75
$ws.Range("M32").Value = -2662.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1621.25
$ws.Range("I45").Value = 1621.25
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1621.25
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1244.25
$ws.Range("N45").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 393.8
$ws.Range("I110").Value = 489
$ws.Range("K110").Value = 489
$ws.Range("M110").Value = 1556

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4304.52
$ws.Range("I122").Value = 4007.9
$ws.Range("K122").Value = 12023.7
$ws.Range("M122").Value = -9573.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1243.75
$ws.Range("I20").Value = 981.44446
$ws.Range("J20").Value = 2030.6666
$ws.Range("K20").Value = 981.44446
$ws.Range("L20").Value = 2030.6666
$ws.Range("M20").Value = -734.44446
$ws.Range("N20").Value = -2524.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1966.6666
$ws.Range("I99").Value = 1900
$ws.Range("K99").Value = 1900
$ws.Range("M99").Value = -402

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2315.7144
$ws.Range("I107").Value = 2492
$ws.Range("J107").Value = 1875
$ws.Range("K107").Value = 2492
$ws.Range("L107").Value = 1875
$ws.Range("M107").Value = -572
$ws.Range("N107").Value = -5715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3506.75
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 76.333336
$ws.Range("J7").Value = 79.333336
$ws.Range("L7").Value = 79.333336
$ws.Range("N7").Value = -305.333336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 450
$ws.Range("I38").Value = 450
$ws.Range("K38").Value = 450
$ws.Range("M38").Value = -73

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H46").Value = 450
$ws.Range("I46").Value = 450
$ws.Range("K46").Value = 450
$ws.Range("M46").Value = -239

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3490.8333
$ws.Range("I58").Value = 1236.25
$ws.Range("K58").Value = 1236.25
$ws.Range("M58").Value = -1033.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4250
$ws.Range("I86").Value = 3500
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 3500
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -2377
$ws.Range("N86").Value = -7246

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 4250
$ws.Range("I89").Value = 3500
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 17500
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -11884
$ws.Range("N89").Value = -36232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2537.375
$ws.Range("I105").Value = 2728.4285
$ws.Range("K105").Value = 2728.4285
$ws.Range("M105").Value = -981.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5567
$ws.Range("I134").Value = 3962
$ws.Range("J134").Value = 6637
$ws.Range("K134").Value = 11886
$ws.Range("L134").Value = 19911
$ws.Range("M134").Value = -9351
$ws.Range("N134").Value = -24981

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3490.8333
$ws.Range("I136").Value = 1236.25
$ws.Range("K136").Value = 3708.75
$ws.Range("M136").Value = -1158.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1039.75
$ws.Range("I140").Value = 1039.75
$ws.Range("K140").Value = 3119.25
$ws.Range("M140").Value = 2060.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = $null
$ws.Range("N2").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 399
$ws.Range("I29").Value = 399
$ws.Range("K29").Value = 399
$ws.Range("M29").Value = -109

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2426.1667
$ws.Range("I102").Value = 1886
$ws.Range("K102").Value = 1886
$ws.Range("M102").Value = -264

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 96600
$ws.Range("J135").Value = 96600
$ws.Range("L135").Value = 96600
$ws.Range("N135").Value = -106740

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 16400
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10334.333
$ws.Range("I7").Value = 3001.5
$ws.Range("J7").Value = 25000
$ws.Range("K7").Value = 3001.5
$ws.Range("L7").Value = 25000
$ws.Range("M7").Value = -2889.5
$ws.Range("N7").Value = -25224

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H28").Value = 16400
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H37").Value = 16400
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1050
$ws.Range("I61").Value = 1050
$ws.Range("K61").Value = 1050
$ws.Range("M61").Value = -848

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1050
$ws.Range("I113").Value = 1050
$ws.Range("K113").Value = 1050
$ws.Range("M113").Value = 1120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 10334.333
$ws.Range("I126").Value = 3001.5
$ws.Range("J126").Value = 25000
$ws.Range("K126").Value = 9004.5
$ws.Range("L126").Value = 75000
$ws.Range("M126").Value = -6534.5
$ws.Range("N126").Value = -79940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 30000000
$ws.Range("J4").Value = 30000000
$ws.Range("L4").Value = 30000000
$ws.Range("N4").Value = -30000226

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 22800
$ws.Range("J21").Value = 22800
$ws.Range("L21").Value = 22800
$ws.Range("N21").Value = -23270

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 3000000
$ws.Range("I29").Value = 3000000
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 3000000
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -2999710
$ws.Range("N29").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H35").Value = 22800
$ws.Range("J35").Value = 22800
$ws.Range("L35").Value = 22800
$ws.Range("N35").Value = -23380

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1039.4286
$ws.Range("I100").Value = 1255.2
$ws.Range("K100").Value = 2510.4
$ws.Range("M100").Value = -1969.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 790.75
$ws.Range("I122").Value = 754.3333
$ws.Range("K122").Value = 2262.9999
$ws.Range("M122").Value = 187.0001000000002
